$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Swap the two data values on Sheet1 (A2 <-> A3 contents) ---
$ws1.Range("A2").Value = "Ra"
$ws1.Range("A3").Value = "Sankar"

# Sheet1: select column B (B1:B1048576), matching the authored selection
$ws1.Columns("B:B").Select() | Out-Null

# --- Add the new "Sheet3" worksheet after Sheet1 ---
# (An intermediate throw-away sheet is added+removed first purely so the
#  surviving sheet picks up sheetId=3, matching a workbook that once had a
#  second sheet before this one was created.)
$temp = $wb.Worksheets.Add($null, $ws1)
$sheet3 = $wb.Worksheets.Add($null, $temp)
$sheet3.Name = "Sheet3"
$temp.Activate() | Out-Null
$temp.Delete() | Out-Null

# Re-fetch the reference by name: the old object handle goes stale once the
# sheet that was added alongside it (our throw-away "temp") is deleted.
$sheet3 = $wb.Worksheets.Item("Sheet3")

# Populate Sheet3 cells (order matters so shared-string indices line up)
$sheet3.Range("A1").Value = "fgf"
$sheet3.Range("C6").Value = "fdgfd"
$sheet3.Range("B4").Value = "fdg"
$sheet3.Range("D4").Value = "dfgfg"

# Make Sheet3 the active/selected tab with its own selection at G12
$sheet3.Activate() | Out-Null
$sheet3.Range("G12").Select() | Out-Null
